$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 3 achievement/gap values
$ws.Range("E3").Value = 4500
$ws.Range("J3").Value = 350

# Row "Valh" (row 37) was removed; rows 38-51 shift up into 37-50.
# Apply the resulting values directly (row 51 ends up empty).

$ws.Range("A37").Value = "xbladze"
$ws.Range("B37").Value = "15/12/2025"
$ws.Range("E37").Value = 109220
$ws.Range("H37").Value = 3850
$ws.Range("I37").Value = 10000
$ws.Range("J37").Value = 0

$ws.Range("A38").Value = "zordan"
$ws.Range("B38").Value = "16/01/2026"
$ws.Range("E38").Value = 67000
$ws.Range("H38").ClearContents()
$ws.Range("I38").Value = 10000
$ws.Range("J38").Value = 8050

$ws.Range("A39").Value = "Amazonie"
$ws.Range("B39").Value = "21/02/2026"
$ws.Range("E39").Value = 73350
$ws.Range("I39").ClearContents()
$ws.Range("J39").Value = 1650

$ws.Range("A40").Value = "terracrom"
$ws.Range("B40").Value = "21/02/2026"
$ws.Range("E40").Value = 20875
$ws.Range("J40").Value = 0

$ws.Range("A41").Value = "Anto"
$ws.Range("B41").Value = "21/02/2026"
$ws.Range("E41").Value = 7650

$ws.Range("A42").Value = "Dasters79"
$ws.Range("B42").Value = "21/02/2026"
$ws.Range("E42").Value = 23855
$ws.Range("J42").Value = 950

$ws.Range("A43").Value = "dibba10"
$ws.Range("B43").Value = "21/02/2026"
$ws.Range("E43").Value = 12460
$ws.Range("J43").Value = 0

$ws.Range("A44").Value = "cucco"
$ws.Range("B44").Value = "21/02/2026"
$ws.Range("E44").Value = 23050
$ws.Range("J44").Value = 1700

$ws.Range("A45").Value = "fede61mito"
$ws.Range("B45").Value = "21/02/2026"
$ws.Range("E45").Value = 400
$ws.Range("J45").Value = 0

$ws.Range("A46").Value = "Xx_Herman_xX"
$ws.Range("B46").Value = "21/02/2026"
$ws.Range("E46").Value = 11735
$ws.Range("J46").Value = 500

$ws.Range("A47").Value = "Artur"
$ws.Range("B47").Value = "21/02/2026"
$ws.Range("E47").Value = 1800
$ws.Range("J47").Value = 0

$ws.Range("A48").Value = "Michele"
$ws.Range("B48").Value = "21/02/2026"
$ws.Range("E48").Value = 450

$ws.Range("A49").Value = "MIRIAM MIRIAM"
$ws.Range("B49").Value = "21/02/2026"
$ws.Range("E49").Value = 19750
$ws.Range("J49").Value = 300

$ws.Range("A50").Value = "SanBucaツ"
$ws.Range("B50").Value = "22/02/2026"
$ws.Range("E50").Value = 43290
$ws.Range("J50").Value = 0

# Old row 51 data is gone entirely now - clear it but keep its row formatting
# (height / thick bottom border) intact.
$ws.Range("A51:W51").Clear()
